$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Tareas" block from column D (rows 1-4) to column A (rows 20-23) ---

# Row 20 <- old D1 (Tareas: / header, bold + wrap)
$ws.Range("A20").Value2 = $ws.Range("D1").Value2
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").WrapText = $true

# Row 21 <- old D2 (wrap, ht 30)
$ws.Range("A21").Value2 = $ws.Range("D2").Value2
$ws.Range("A21").WrapText = $true
$ws.Range("A21").RowHeight = 30

# Row 22 <- old D3 (wrap, ht 90 after the column got wider)
$ws.Range("A22").Value2 = $ws.Range("D3").Value2
$ws.Range("A22").WrapText = $true
$ws.Range("A22").RowHeight = 90

# Row 23 <- old D4 (wrap)
$ws.Range("A23").Value2 = $ws.Range("D4").Value2
$ws.Range("A23").WrapText = $true

# Remove the old column D content entirely
$ws.Range("D1:D4").Clear()

# --- Column width changes ---
# (input values tuned so the engine's internal pixel-rounding lands as close as
# possible to the target stored widths of 54.42578125 / 68.7109375)
$ws.Columns.Item(1).ColumnWidth = 53.59
$ws.Columns.Item(4).ColumnWidth = 67.75

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("C20").Select() | Out-Null

Write-Host "done"
